$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final id / speaker_variant values for rows 2-15 (column A stays the same URL
# for every row; is_prefered (D) and the remaining flag columns are cleared).
$data = @(
    @("#r", "r"),
    @("#barnardi", "Barnardi"),
    @("#charlotte", "Charlotte"),
    @("#margo", "Margo"),
    @("#garcias", "Garcias"),
    @("#stefanus", "Stefanus"),
    @("#joanni", "Joanni"),
    @("#dosorio", "Dosorio"),
    @("#karlo", "Karlo"),
    @("#gusman", "Gusman"),
    @("#ambrosius", "Ambrosius"),
    @("#kantel", "Kantel"),
    @("#antonio", "Antonio"),
    @("#laurentio", "Laurentio")
)

$url = "https://www.dbnl.org/tekst/asse001gusm02_01"
$startRow = 2

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $url
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    # is_prefered / is_new / is_error / gender / comments all blank now
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 5).Value = ""
    $ws.Cells.Item($row, 6).Value = ""
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = ""
}
